$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 8 data: "NILatticeKd" ---
$ws.Range("A8").Value = "NILatticeKd"
$ws.Range("B8").Value = 820
$ws.Range("C8").Value = 1702
$ws.Range("D8").Value = 1707
$ws.Range("E8").Value = 1819
$ws.Range("F8").Value = 2160
$ws.Range("G8").Value = 2560
$ws.Range("H8").Value = 349.1
$ws.Range("I8").Value = "[]"
$ws.Range("J8").Value = "[1 1]"
$ws.Range("K8").Value = "[100 100]"

# --- Column width adjustments (closest achievable given engine's char-width rounding) ---
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(6).ColumnWidth = 10.333333333333334
$ws.Columns.Item(7).ColumnWidth = 10.5
$ws.Columns.Item(8).ColumnWidth = 5.5
$ws.Columns.Item(9).ColumnWidth = 16
$ws.Columns.Item(10).ColumnWidth = 18.666666666666668
$ws.Columns.Item(11).ColumnWidth = 16

# --- Update selection to A8 ---
$ws.Range("A8").Select()
